# Weekly update: insert a new "Agrícola del Norte S.A. de Arica - Zanahoria"
# record at row 178, pushing the existing rows (178-218) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row before the current row 178; this shifts rows 178..218
# down to 179..219 and carries the date number format (style id 2) on
# column D along with it.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with this week's data.
$ws.Cells.Item(178, 1).Value = 1
$ws.Cells.Item(178, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(178, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(178, 4).Value = 44543
$ws.Cells.Item(178, 5).Value = 15
$ws.Cells.Item(178, 6).Value = 100114013
$ws.Cells.Item(178, 7).Value = "Zanahoria"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 90
$ws.Cells.Item(178, 11).Value = 18000
$ws.Cells.Item(178, 12).Value = 19000
$ws.Cells.Item(178, 13).Value = 18500
$ws.Cells.Item(178, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(178, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(178, 16).Value = 740
$ws.Cells.Item(178, 17).Value = 25
$ws.Cells.Item(178, 18).Value = "Hortaliza"
